$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '26.133.59'
Set-TextValue $ws 'D3' '1.670.04'
Set-TextValue $ws 'E3' '  -1.35%  '
Set-TextValue $ws 'E4' '  -0.74%  '
Set-TextValue $ws 'D5' '210.74'
Set-TextValue $ws 'E5' '  -3.76%  '
Set-TextValue $ws 'D6' '0.5209'
Set-TextValue $ws 'E6' '  -5.02%  '
Set-TextValue $ws 'E7' '  -0.71%  '
Set-TextValue $ws 'D8' '0.2639'
Set-TextValue $ws 'E8' '  -3.39%  '
Set-TextValue $ws 'D9' '0.06239'
Set-TextValue $ws 'E9' '  -3.58%  '
Set-TextValue $ws 'D10' '21.16'
Set-TextValue $ws 'E10' '  -3.92%  '
Set-TextValue $ws 'D11' '0.07509'
Set-TextValue $ws 'E11' '  -2.13%  '
Set-TextValue $ws 'D12' '1.652.26'
Set-TextValue $ws 'E12' '  -2.82%  '
Set-TextValue $ws 'D13' '4.437'
Set-TextValue $ws 'E13' '  -2.55%  '
Set-TextValue $ws 'D14' '0.5589'
Set-TextValue $ws 'E14' '  -4.40%  '
Set-TextValue $ws 'D15' '66.23'
Set-TextValue $ws 'E15' '  +1.29%  '
Set-TextValue $ws 'D16' '0.000007953'
Set-TextValue $ws 'E16' '  -5.26%  '
Set-TextValue $ws 'D17' '26.178.65'
Set-TextValue $ws 'E17' '  -1.00%  '
Set-TextValue $ws 'E18' '  -0.77%  '
Set-TextValue $ws 'D19' '4.792'
Set-TextValue $ws 'E19' '  -3.23%  '
Set-TextValue $ws 'D20' '186.86'
Set-TextValue $ws 'E20' '  -2.73%  '
Set-TextValue $ws 'E21' '  -5.67%  '
Set-TextValue $ws 'D22' '6.176'
Set-TextValue $ws 'E22' '  -1.33%  '
Set-TextValue $ws 'E23' '  -0.72%  '
Set-TextValue $ws 'D24' '147.68'
Set-TextValue $ws 'E24' '  -1.28%  '
Set-TextValue $ws 'E25' '  -6.36%  '
Set-TextValue $ws 'D26' '7.584'
Set-TextValue $ws 'E26' '  -4.00%  '
Set-TextValue $ws 'D27' '15.87'
Set-TextValue $ws 'E27' '  +0.76%  '
Set-TextValue $ws 'D28' '0.06195'
Set-TextValue $ws 'E28' '  -2.03%  '
Set-TextValue $ws 'D29' '1.355'
Set-TextValue $ws 'E29' '  -3.03%  '
Set-TextValue $ws 'D30' '1.279'
Set-TextValue $ws 'E30' '  -3.93%  '
Set-TextValue $ws 'D31' '3.475'
Set-TextValue $ws 'E31' '  -3.62%  '
Set-TextValue $ws 'D32' '3.424'
Set-TextValue $ws 'E32' '  -4.91%  '
Set-TextValue $ws 'D33' '1.606'
Set-TextValue $ws 'E33' '  -4.68%  '
Set-TextValue $ws 'D34' '0.9915'
Set-TextValue $ws 'E34' '  -5.09%  '
Set-TextValue $ws 'D35' '0.6030'
Set-TextValue $ws 'E35' '  -2.05%  '
Set-TextValue $ws 'D36' '2.405'
Set-TextValue $ws 'E36' '  -0.20%  '
Set-TextValue $ws 'D37' '2.705'
Set-TextValue $ws 'E37' '  -0.15%  '
Set-TextValue $ws 'D38' '6.128'
Set-TextValue $ws 'E38' '  -1.43%  '
Set-TextValue $ws 'D39' '0.01608'
Set-TextValue $ws 'E39' '  -1.70%  '
Set-TextValue $ws 'D40' '0.8675'
Set-TextValue $ws 'E40' '  -2.05%  '
Set-TextValue $ws 'D41' '1.072.23'
Set-TextValue $ws 'E41' '  -4.25%  '
Set-TextValue $ws 'D43' '99.51'
Set-TextValue $ws 'E43' '  -2.53%  '
Set-TextValue $ws 'D44' '1.818.83'
Set-TextValue $ws 'D45' '0.00000000107'
Set-TextValue $ws 'E45' '  +0.52%  '
Set-TextValue $ws 'D46' '55.97'
Set-TextValue $ws 'E46' '  -2.74%  '
Set-TextValue $ws 'E47' '  -0.08%  '
Set-TextValue $ws 'E48' '  -0.59%  '
Set-TextValue $ws 'D49' '7.917'
Set-TextValue $ws 'E49' '  -3.77%  '
Set-TextValue $ws 'E50' '  -1.19%  '
Set-TextValue $ws 'D51' '5.951'
Set-TextValue $ws 'E51' '  -2.48%  '
